# Adds an estimated feed-mass-eaten-per-pig column to the "individual" sheet,
# adds a purchase_date column, renames "misc" -> "feed_per_pig" on the monthly
# sheets, and brings the "2022" sheet's layout + data in line with "2021".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "individual"
# ---------------------------------------------------------------------------
$ind = $wb.Worksheets.Item("individual")

# New headers: J1 = feed_eaten, K1 = purchase_date (styled like the other
# yellow header cells), L1 picks up a stray numeric marker (9) that used to
# live in J1.
$ind.Range("J1").Value = "feed_eaten"
$ind.Range("K1").Value = "purchase_date"
$ind.Range("J1:K1").Interior.Color = 65535
$ind.Range("L1").Value = 9

# feed_eaten (J) values, estimated feed mass eaten per pig.
$ind.Range("J2").Value = 115
$ind.Range("J3").Value = 110
$ind.Range("J4").Value = 110
$ind.Range("J5").Value = 108
$ind.Range("J6").Value = 108
$ind.Range("J7").Value = 0
$ind.Range("J8").Value = 0
$ind.Range("J9").Value = 0
$ind.Range("J10").Value = 0

# purchase_date (K) values - same dates as when each batch's birth_date was
# recorded for rows 2-6, and the later June purchase date for rows 7-10.
$ind.Range("K2").Value = 44141
$ind.Range("K2").NumberFormat = "dd/mm/yy"
$ind.Range("K3").Value = 44275
$ind.Range("K3").NumberFormat = "yyyy\-mm\-dd"
$ind.Range("K4").Value = 44275
$ind.Range("K4").NumberFormat = "yyyy\-mm\-dd"
$ind.Range("K5").Value = 44352
$ind.Range("K5").NumberFormat = "yyyy\-mm\-dd"
$ind.Range("K6").Value = 44352
$ind.Range("K6").NumberFormat = "yyyy\-mm\-dd"
$ind.Range("K7").Value = 44427
$ind.Range("K7").NumberFormat = "yyyy\-mm\-dd"
$ind.Range("K8").Value = 44427
$ind.Range("K8").NumberFormat = "yyyy\-mm\-dd"
$ind.Range("K9").Value = 44427
$ind.Range("K9").NumberFormat = "yyyy\-mm\-dd"
$ind.Range("K10").Value = 44427
$ind.Range("K10").NumberFormat = "yyyy\-mm\-dd"

# meds values for row 8 and row 10 were swapped during review.
$ind.Range("I8").Value = 10
$ind.Range("I10").Value = 40

# Touch G9/G10 so they exist as real (blank) cells, matching column G's
# formatting elsewhere in the table.
$ind.Range("G9").NumberFormat = "General"
$ind.Range("G10").NumberFormat = "General"

# A new trailing blank row (11) with the same per-column formatting as the
# rows above it (date formats under the date columns).
$ind.Range("A11").NumberFormat = "General"
$ind.Range("B11").NumberFormat = "yyyy-mm-dd"
$ind.Range("D11").NumberFormat = "yyyy-mm-dd"
$ind.Range("H11").NumberFormat = "General"
$ind.Range("I11").NumberFormat = "General"
$ind.Range("K11").NumberFormat = "yyyy-mm-dd"

$ind.Range("E1").Select()

# ---------------------------------------------------------------------------
# Sheet "2021"
# ---------------------------------------------------------------------------
$y21 = $wb.Worksheets.Item("2021")

$y21.Range("E1").Value = "feed_per_pig"

$y21.Range("E6").Value = 50
$y21.Range("E7").Value = 33
$y21.Range("E8").Value = 33
$y21.Range("F8").NumberFormat = "General"
$y21.Range("F8").Value = 71

$y21.Range("C9").Value = 250
$y21.Range("D9").Value = 1606
$y21.Range("E9").Value = 42
$y21.Range("F9").NumberFormat = "General"
$y21.Range("F9").Value = 71

$y21.Range("C10").Value = 150
$y21.Range("D10").Value = 937
$y21.Range("E10").Value = 38

$y21.Range("E11").ClearContents()
$y21.Range("E12").ClearContents()
$y21.Range("E13").ClearContents()

$y21.Range("A9").Select()

# ---------------------------------------------------------------------------
# Sheet "2022" - rebuild to mirror "2021"'s layout (month rows instead of
# month columns) and copy this year's figures across.
# ---------------------------------------------------------------------------
$y22 = $wb.Worksheets.Item("2022")

# Clear the old transposed table first.
$y22.Range("A1:M8").ClearContents()

$y22.Range("A1").Value = "month"
$y22.Range("B1").Value = "population"
$y22.Range("C1").Value = "feed_mass"
$y22.Range("D1").Value = "feed_price"
$y22.Range("E1").Value = "feed_per_pig"
$y22.Range("F1").Value = "average_age"

$months = @("jan","feb","mar","apr","may","jun","jul","aug","sep","oct","nov","dec")
for ($i = 0; $i -lt 12; $i++) {
    $y22.Cells.Item($i + 2, 1).Value = $months[$i]
}

$y22.Range("B2").Value = 0
$y22.Range("C2").Value = 0
$y22.Range("D2").Value = 0
$y22.Range("E2").Value = 0
$y22.Range("F2").Value = 0

$y22.Range("B3").Value = 0
$y22.Range("C3").Value = 0
$y22.Range("D3").Value = 0
$y22.Range("E3").Value = 0
$y22.Range("F3").Value = 0

$y22.Range("B4").Value = 0
$y22.Range("C4").Value = 0
$y22.Range("D4").Value = 0
$y22.Range("E4").Value = 0
$y22.Range("F4").Value = 0

$y22.Range("B5").Value = 0
$y22.Range("C5").Value = 0
$y22.Range("D5").Value = 0
$y22.Range("E5").Value = 0
$y22.Range("F5").Value = 0

$y22.Range("B6").Value = 3
$y22.Range("C6").Value = 150
$y22.Range("D6").Value = 960
$y22.Range("E6").Value = 50
$y22.Range("F6").Value = 96

$y22.Range("B7").Value = 3
$y22.Range("C7").Value = 100
$y22.Range("D7").Value = 624
$y22.Range("E7").Value = 33
$y22.Range("F7").Value = 77

$y22.Range("B8").Value = 6
$y22.Range("C8").Value = 200
$y22.Range("D8").Value = 1292
$y22.Range("E8").Value = 33
$y22.Range("F8").Value = 71

$y22.Range("B9").Value = 6
$y22.Range("C9").Value = 0
$y22.Range("D9").Value = 0
$y22.Range("F9").Value = 71

$y22.Range("B10").Value = 4
$y22.Range("C10").Value = 50
$y22.Range("D10").Value = 312
$y22.Range("E10").Value = 10
$y22.Range("F10").Value = 56

$y22.Range("B11").Value = 4
$y22.Range("C11").Value = 0
$y22.Range("D11").Value = 0

$y22.Range("B12").Value = 0
$y22.Range("C12").Value = 0
$y22.Range("D12").Value = 0

$y22.Range("B13").Value = 0
$y22.Range("C13").Value = 0
$y22.Range("D13").Value = 0

$y22.Range("A1").Select()
